$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Handgin" -> "Handgun" typo and re-order the underlying shared
# strings so that "Restrictions on Handgun Ownership" becomes the newest
# (last) shared string while the other two category labels keep their text
# but pick up new shared-string slots.
$ws.Range("B11").Value = "Restrictions on Handgun Ownership"
$ws.Range("B12").Value = "Restrictions on Carrying Long Guns"
$ws.Range("B13").Value = "Restrictions on Carrying Handguns"

# Move the selection/active cell to B12, matching the saved view state.
$ws.Range("B12").Select()
